$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1), matching the style of the existing
# header cells (copy format from G1 so it gets the same bold/border/
# alignment style used by the other headers).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# "Save" flag values for the two data rows.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
